$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" value (row 8, col B) to the new
#    generation timestamp.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-22T16:25:12+00:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: the "Mapping: RIM Mapping" column and the
#    "Mapping: Spécification métier vers l'extension ROR LocationStatus"
#    column (AK and AL) were swapped - header, every data row, and the
#    column widths that were sized to fit each column's content.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 6; $r++) {
    $akCell = $ws.Range("AK$r")
    $alCell = $ws.Range("AL$r")
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    if ($akVal -ne $alVal) {
        $akCell.Value2 = $alVal
        $alCell.Value2 = $akVal
    }
}

# Column widths follow the content: AK used to be the narrow column
# (~24.98) and AL the wide one (~70.57); after the swap AK is wide and
# AL is narrow.
$ws.Range("AK1").ColumnWidth = 69.66666666666667
$ws.Range("AL1").ColumnWidth = 24.166666666666668
